$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 379, shifting rows 379:404 down to 380:405
$ws.Rows.Item(379).Insert()

# Populate the new row 379 with the new data
$ws.Cells.Item(379, 1).Value = 7
$ws.Cells.Item(379, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(379, 3).Value = "Ñuble"
$ws.Cells.Item(379, 4).Value = 44826
$ws.Cells.Item(379, 5).Value = 16
$ws.Cells.Item(379, 6).Value = 100114001
$ws.Cells.Item(379, 7).Value = "Papa"
$ws.Cells.Item(379, 8).Value = "Patagonia"
$ws.Cells.Item(379, 9).Value = "1a (guarda)"
$ws.Cells.Item(379, 10).Value = 160
$ws.Cells.Item(379, 11).Value = 7000
$ws.Cells.Item(379, 12).Value = 7500
$ws.Cells.Item(379, 13).Value = 7250
$ws.Cells.Item(379, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(379, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(379, 16).Value = 290
$ws.Cells.Item(379, 17).Value = 25
$ws.Cells.Item(379, 18).Value = "Hortaliza"
